# Chatbot test content updates
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "greetings.hi" query gained two more trigger phrases
$ws.Range("B2").Value = "hi; hello; how are you doing?; how is it going?; Hi?; Hello?"

# "ask.what_is_ml" query gained an additional phrasing
$ws.Range("B6").Value = "what is machine learning?; what is ML?;"

# Column B grew wider text, so re-fit its width to the new content
$ws.Columns("B").AutoFit()

# Leave the cursor on B7, matching where editing left off
$ws.Range("B7").Select()
